$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name "BrassA-HW40.xpc" -> "BrassA")
$ws.Name = "BrassA"

# Correct tiny floating point rounding differences in existing cells
$ws.Range("D13").Value = 0.9926792305866939
$ws.Range("H13").Value = 0.9926792305866939
$ws.Range("N13").Value = 0.9925484321381071
$ws.Range("F15").Value = 0.9928943056861634
$ws.Range("J15").Value = 0.9939504500144281
$ws.Range("L15").Value = 0.9963040968158872
$ws.Range("O15").Value = 0.9926730406175921

# Append a new data row (row 16) for the "HexGrid-60degTilt5degRes" case
# Copy the formatting from the row above (bold/bordered/centered "index" style)
# onto the new index cell before writing its value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9923753196759012
$ws.Range("D16").Value = 0.9879300160609247
$ws.Range("E16").Value = 0.995896977005419
$ws.Range("F16").Value = 0.9946912582886986
$ws.Range("G16").Value = 0.9923753196759012
$ws.Range("H16").Value = 0.9879300160609247
$ws.Range("I16").Value = 0.9921929612759364
$ws.Range("J16").Value = 0.9948219900607478
$ws.Range("K16").Value = 0.9946551405311592
$ws.Range("L16").Value = 0.9824036829982525
$ws.Range("M16").Value = 0.9923753196759012
$ws.Range("N16").Value = 0.9919134965331718
$ws.Range("O16").Value = 0.9927233927577359
$ws.Range("P16").Value = 0.9918709182371299
